# Update the multiplication problems/answers in the document to the newly
# generated set, per the commit "Update master to output generated at 9a8706d".

$d = $word.ActiveDocument

$replacements = @(
    @{ old = "75×83=6225"; new = "62×61=3782" },
    @{ old = "82×16=1312"; new = "54×68=3672" },
    @{ old = "53×85=4505"; new = "64×35=2240" },
    @{ old = "75×56=4200"; new = "70×56=3920" },
    @{ old = "55×46=2530"; new = "29×30=870" },
    @{ old = "38×64=2432"; new = "98×66=6468" },
    @{ old = "71×20=1420"; new = "45×41=1845" },
    @{ old = "36×46=1656"; new = "41×98=4018" },
    @{ old = "77×94=7238"; new = "22×93=2046" },
    @{ old = "62×27=1674"; new = "81×79=6399" },
    @{ old = "28×61=1708"; new = "71×63=4473" },
    @{ old = "46×51=2346"; new = "99×19=1881" },
    @{ old = "73×75=5475"; new = "25×17=425" },
    @{ old = "94×24=2256"; new = "83×71=5893" },
    @{ old = "23×87=2001"; new = "44×85=3740" },
    @{ old = "84×22=1848"; new = "25×34=850" },
    @{ old = "43×25=1075"; new = "90×96=8640" },
    @{ old = "14×81=1134"; new = "47×18=846" },
    @{ old = "87×89=7743"; new = "43×60=2580" },
    @{ old = "67×75=5025"; new = "77×60=4620" },
    @{ old = "24×31=744";  new = "53×12=636" },
    @{ old = "63×67=4221"; new = "74×93=6882" },
    @{ old = "93×27=2511"; new = "98×82=8036" },
    @{ old = "27×48=1296"; new = "72×98=7056" },
    @{ old = "45×44=1980"; new = "80×98=7840" }
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.old, $true, $true, $false, $false, $false, $true, 1, $false, $r.new, 2)
}

$d.Save()
